$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <-> Row 3: course/student-count swap; times become "HH:MM" text.
# Venue/capacity (F/G) stay attached to the row, not the course.
$ws.Range("B2").Value = "MAT141"
$ws.Range("C2").Value = "08:00"
$ws.Range("D2").Value = "11:00"
$ws.Range("E2").Value = 30

$ws.Range("B3").Value = "SOC205"
$ws.Range("C3").Value = "08:30"
$ws.Range("D3").Value = "11:30"
$ws.Range("E3").Value = 40

# Row 4: PHY201 - times become "HH:MM" text only
$ws.Range("C4").Value = "09:00"
$ws.Range("D4").Value = "12:00"

# Row 5: ECON202 - times become "HH:MM" text only
$ws.Range("C5").Value = "09:30"
$ws.Range("D5").Value = "12:30"

# Row 6: CSC103 - times become "HH:MM" text only
$ws.Range("C6").Value = "10:00"
$ws.Range("D6").Value = "13:00"

# Row 7: ENG220 - times become "HH:MM" text only
$ws.Range("C7").Value = "10:30"
$ws.Range("D7").Value = "13:30"

# Row 8: STA121 - times become "HH:MM" text only; venue_capacity "null" -> 0
$ws.Range("C8").Value = "11:00"
$ws.Range("D8").Value = "14:00"
$ws.Range("G8").Value = 0

# Row 9: GEO111 - times become "HH:MM" text only
$ws.Range("C9").Value = "11:30"
$ws.Range("D9").Value = "14:30"

# Row 10: HIS101 - times become "HH:MM" text only
$ws.Range("C10").Value = "12:00"
$ws.Range("D10").Value = "15:00"

# Row 11: MUS120 - times become "HH:MM" text only
$ws.Range("C11").Value = "12:30"
$ws.Range("D11").Value = "15:30"

# Row 12 <-> Row 13: course swap, student counts swap, times become "HH:MM" text
$ws.Range("B12").Value = "BIO101"
$ws.Range("C12").Value = "13:00"
$ws.Range("D12").Value = "16:00"
$ws.Range("E12").Value = 45

$ws.Range("B13").Value = "PSY101"
$ws.Range("C13").Value = "13:30"
$ws.Range("D13").Value = "16:30"
$ws.Range("E13").Value = 65

# Row 14: CSC103 - times become "HH:MM" text; venue/capacity swap with row 15
$ws.Range("C14").Value = "14:00"
$ws.Range("D14").Value = "17:00"
$ws.Range("F14").Value = "CBN"
$ws.Range("G14").Value = 500

# Row 15: CHEM301 - times become "HH:MM" text; venue/capacity swap with row 14
$ws.Range("C15").Value = "14:30"
$ws.Range("D15").Value = "17:30"
$ws.Range("F15").Value = "FLT"
$ws.Range("G15").Value = 100

# Row 16: ART150 - times become "HH:MM" text; venue/capacity updated
$ws.Range("C16").Value = "15:00"
$ws.Range("D16").Value = "18:00"
$ws.Range("F16").Value = "SLT"
$ws.Range("G16").Value = 100
